# Update the "想去人数" (want-to-go count) figures in both the "展览"
# and "全部类型" worksheets to reflect newly scraped values.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 6765
$ws1.Range("F7").Value = 10
$ws1.Range("F10").Value = 6345
$ws1.Range("F15").Value = 105
$ws1.Range("F25").Value = 158

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 6765
$ws4.Range("F7").Value = 10
$ws4.Range("F10").Value = 6345
$ws4.Range("F15").Value = 105
$ws4.Range("F26").Value = 158
